$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New API rows describing patient registration / update / delete endpoints.
# Cell values are set in the same order the author typed them so the
# shared-strings table is interned in the matching sequence.
$ws.Range("B4").Value = "POST"
$ws.Range("D4").Value = "환자 추가"

$ws.Range("B6").Value = "DELETE"

$ws.Range("B5").Value = "PUT"
$ws.Range("D5").Value = "환자 정보 수정"

$ws.Range("D6").Value = "환자 삭제"

# Index numbers and URLs (URLs reuse strings already present in the sheet).
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5

$ws.Range("C4").Value = "/api/patients "
$ws.Range("C5").Value = "/api/patients/<int:pk> "
$ws.Range("C6").Value = "/api/patients/<int:pk> "

# Rows 4-23, columns A-C pick up the centered style used by the existing
# data rows (2-3); column D keeps its current left/general-aligned style.
$ws.Range("A4:C23").HorizontalAlignment = -4108
